# Add a new menu row (id, name, price, branch, category, description) on
# the "menu" sheet. Row 12 was an empty placeholder row (only cell A12 had
# a leftover red "missing value" style); fill it in with a new "Coke"
# entry for the NTU branch's Drink category, sold at half price.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the placeholder row's leftover formatting (the red font style on
# A12) so the new row matches the plain look of every other data row.
$ws.Range("A12:F12").ClearFormats()

# The price column stores values as plain text elsewhere in this sheet
# (e.g. "3.2", "9.9", ...). Mark C12 as Text before writing "2.0" so Excel
# keeps it as the literal string "2.0" instead of auto-converting it to
# the number 2.
$ws.Range("C12").NumberFormat = "@"

$ws.Range("A12").Value = "ed718b0b-88f1-44cb-b4db-493fcf3ff473"
$ws.Range("B12").Value = "Coke"
$ws.Range("C12").Value = "2.0"
$ws.Range("D12").Value = "NTU"
$ws.Range("E12").Value = "Drink"
$ws.Range("F12").Value = "Half Price!"

# Drop the temporary Text number-format override so C12 ends up with the
# same default (unstyled) formatting as its neighbours.
$ws.Range("C12").ClearFormats()
